{"js": "// Replace the title prefix \"Video #17 - \" with \"# - \" while keeping the\n// rest of the title (\"Users logged in for 5 consecutive days\") intact.\nconst results = context.document.body.search(\"Video #17 - \", { matchCase: true });\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length > 0) {\n  results.items[0].insertText(\"# - \", Word.InsertLocation.replace);\n  await context.sync();\n}\n", "ps1": "# Replace the title prefix \"Video #17 - \" with \"# - \" while keeping the\n# rest of the title (\"Users logged in for 5 consecutive days\") intact.\n$d = $word.ActiveDocument\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"Video #17 - \"\n$find.Replacement.Text = \"# - \"\n$find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2)\n"}
